# Added the logic of Hashtable in Dataprovider
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testdata")

# Rename the data-provider key cells to the new test-case style naming
$ws.Range("A1").Value = "test01_ValidCreateCustomerAPI"
$ws.Range("A7").Value = "test02_InvalidCreateCustomerAPI"

# Column A needs to be a bit wider to fit the longer names
$ws.Columns.Item(1).ColumnWidth = 28.166666666666668

# Move the active selection to A16 (next empty row below the data)
$ws.Range("A16").Select()
